$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.715.90"
$ws.Range("E2").Value = "  +2.53%  "
$ws.Range("D3").Value = "2.413.67"
$ws.Range("E3").Value = "  +2.92%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'554.38"
$ws.Range("E5").Value = "  +2.27%  "
$ws.Range("E6").Value = "  +1.42%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +1.43%  "
$ws.Range("E9").Value = "  +5.27%  "
$ws.Range("E10").Value = "  +2.34%  "
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("E12").Value = "  -2.01%  "
$ws.Range("E13").Value = "  +3.40%  "
$ws.Range("D14").Value = "2.838.73"
$ws.Range("E14").Value = "  +2.85%  "
$ws.Range("D15").Value = "59.575.86"
$ws.Range("E15").Value = "  +2.46%  "
$ws.Range("D16").Value = "'0.0000139"
$ws.Range("E16").Value = "  +4.44%  "
$ws.Range("D17").Value = "2.403.24"
$ws.Range("E17").Value = "  +3.66%  "
$ws.Range("D18").Value = "'11.31"
$ws.Range("E18").Value = "  +5.44%  "
$ws.Range("D19").Value = "'4.45"
$ws.Range("E19").Value = "  +4.50%  "
$ws.Range("D20").Value = "'336.12"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").Value = "'6.97"
$ws.Range("E21").Value = "  +4.61%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "'64.62"
$ws.Range("E23").Value = "  +2.96%  "
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("D25").Value = "'8.50"
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -2.36%  "
$ws.Range("D28").Value = "0.0₃0785"
$ws.Range("E28").Value = "  +6.59%  "
$ws.Range("E29").Value = "  +2.56%  "
$ws.Range("D30").Value = "'170.55"
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("E31").Value = "  +2.33%  "
$ws.Range("D32").Value = "'18.71"
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("D33").Value = "'1.02"
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  +5.39%  "
$ws.Range("D36").Value = "'4.30"
$ws.Range("E36").Value = "  +1.19%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").Value = "'40.13"
$ws.Range("E39").Value = "  +2.61%  "
$ws.Range("D40").Value = "'0.419"
$ws.Range("E40").Value = "  +11.55%  "
$ws.Range("D41").Value = "'305.62"
$ws.Range("E41").Value = "  +6.45%  "
$ws.Range("E42").Value = "  +2.91%  "
$ws.Range("D43").Value = "'142.28"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("E45").Value = "  +4.61%  "
$ws.Range("D46").Value = "'0.572"
$ws.Range("E46").Value = "  +1.32%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'19.08"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("B48").Value = "Polygon"
$ws.Range("C48").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D48").Value = "'0.406"
$ws.Range("E48").Value = "  +6.45%  "
$ws.Range("E49").Value = "  +3.49%  "
$ws.Range("D50").Value = "'11.05"
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("D51").Value = "'1.61"
$ws.Range("E51").Value = "  +5.29%  "
